$d = $word.ActiveDocument

# Each pair (old, new) below corresponds to one cell in the worksheet table.
# Order matters: "998÷4=" -> "996÷5=" must run before "291÷3=" -> "998÷4="
# so the newly written "998÷4=" text is not re-matched by the later search.
$replacements = @(
    ,@("138÷9=", "991÷2=")
    ,@("470÷4=", "884÷5=")
    ,@("948÷6=", "334÷9=")
    ,@("595÷7=", "952÷6=")
    ,@("158÷7=", "948÷5=")
    ,@("989÷2=", "397÷4=")
    ,@("545÷8=", "489÷8=")
    ,@("146÷5=", "317÷7=")
    ,@("977÷2=", "210÷3=")
    ,@("200÷3=", "319÷2=")
    ,@("219÷3=", "314÷2=")
    ,@("661÷6=", "584÷5=")
    ,@("998÷4=", "996÷5=")
    ,@("291÷3=", "998÷4=")
    ,@("689÷3=", "504÷9=")
    ,@("770÷2=", "778÷4=")
    ,@("382÷9=", "988÷2=")
    ,@("912÷5=", "174÷8=")
    ,@("302÷7=", "247÷9=")
    ,@("726÷5=", "949÷4=")
    ,@("425÷8=", "138÷6=")
    ,@("659÷9=", "608÷2=")
    ,@("910÷2=", "821÷2=")
    ,@("144÷5=", "214÷2=")
    ,@("932÷9=", "839÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute could not find text: $old"
    }
}

